$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the "Kla Corporation" row's Security Currency (D2) was wrongly
# recorded as USD; correct it to CAD.
$ws.Range("D2").Value = "CAD"

# Leave the cell selection where the author left it when saving (D3).
$ws.Range("D3").Select()
